$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.13"
$ws.Range("E2").Value = "'-0.36%"

$ws.Range("D3").Value = "'40.72"
$ws.Range("E3").Value = "'3.85%"

$ws.Range("D4").Value = "'5.107"
$ws.Range("E4").Value = "'1.79%"

$ws.Range("D5").Value = "'0.07602"
$ws.Range("E5").Value = "'-1.69%"

$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.618"
$ws.Range("E6").Value = "'3.10%"

$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").Value = "'2.459"
$ws.Range("E7").Value = "'-3.91%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9054"
$ws.Range("E8").Value = "'-1.32%"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1015"
$ws.Range("E9").Value = "'-0.41%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1753"
$ws.Range("E10").Value = "'1.55%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.09074"
$ws.Range("E11").Value = "'0.08%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04275"
$ws.Range("E12").Value = "'-4.51%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1056"
$ws.Range("E13").Value = "'-0.22%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001249"
$ws.Range("E14").Value = "'-2.47%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005869"
$ws.Range("E15").Value = "'3.87%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.350"
$ws.Range("E16").Value = "'-0.39%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'4.269"
$ws.Range("E17").Value = "'-0.45%"

$ws.Range("D19").Value = "'6.592"
$ws.Range("E19").Value = "'-6.42%"

$ws.Range("E20").Value = "'-0.53%"

$ws.Range("D21").Value = "'0.2726"
$ws.Range("E21").Value = "'-1.82%"

$ws.Range("D22").Value = "'0.04182"
$ws.Range("E22").Value = "'1.07%"

$ws.Range("D23").Value = "'0.001229"
$ws.Range("E23").Value = "'2.44%"

$ws.Range("E24").Value = "'-0.59%"

$ws.Range("E25").Value = "'6.30%"

$ws.Range("D26").Value = "'0.0003011"
$ws.Range("E26").Value = "'0.66%"

$ws.Range("D38").Value = "'0.02372"
$ws.Range("E38").Value = "'0.66%"

$ws.Range("D39").Value = "'0.05144"
$ws.Range("E39").Value = "'0.28%"

$ws.Range("D40").Value = "'0.007774"
$ws.Range("E40").Value = "'-2.64%"

$ws.Range("D41").Value = "'0.1296"
$ws.Range("E41").Value = "'-2.21%"

$ws.Range("E42").Value = "'-4.38%"

$ws.Range("D43").Value = "'0.001920"
$ws.Range("E43").Value = "'-3.34%"

$ws.Range("D44").Value = "'0.008448"
$ws.Range("E44").Value = "'5.41%"

$ws.Range("D45").Value = "'0.3325"
$ws.Range("E45").Value = "'-0.06%"

$ws.Range("D46").Value = "'0.00006367"
$ws.Range("E46").Value = "'-4.85%"

$ws.Range("E47").Value = "'-0.23%"

$ws.Range("D48").Value = "'0.004405"
$ws.Range("E48").Value = "'7.01%"

$ws.Range("D49").Value = "'0.03154"
$ws.Range("E49").Value = "'829.38%"

$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.23%"

$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.23%"
